$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 9.156959333333335
$ws.Range("H2").Value = 27.470878
$ws.Range("I2").Value = 0.969469463764299
$ws.Range("J2").Value = 0.9694694637642989
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.553279333333334
$ws.Range("N2").Value = 7.659838000000001
$ws.Range("O2").Value = 0.1645043904057808
$ws.Range("P2").Value = 0.1645043904057808
$ws.Range("Q2").Value = 23.38027502197378
$ws.Range("R2").Value = 210.422475197764
$ws.Range("S2").Value = 0.1594819831535652
$ws.Range("T2").Value = 0.1594819831535652

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 9.156959333333335
$ws.Range("H3").Value = 27.470878
$ws.Range("I3").Value = 0.969469463764299
$ws.Range("J3").Value = 0.9694694637642989
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 8.058662
$ws.Range("N3").Value = 24.175986
$ws.Range("O3").Value = 0.5192088709172035
$ws.Range("P3").Value = 0.5192088709172035
$ws.Range("Q3").Value = 73.79284021507868
$ws.Range("R3").Value = 664.1355619357081
$ws.Range("S3").Value = 0.5033571456697684
$ws.Range("T3").Value = 0.5033571456697683

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 9.156959333333335
$ws.Range("H4").Value = 27.470878
$ws.Range("I4").Value = 0.969469463764299
$ws.Range("J4").Value = 0.9694694637642989
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.909099333333334
$ws.Range("N4").Value = 14.727298
$ws.Range("O4").Value = 0.3162867386770157
$ws.Range("P4").Value = 0.3162867386770157
$ws.Range("Q4").Value = 44.95242295862712
$ws.Range("R4").Value = 404.5718066276441
$ws.Range("S4").Value = 0.3066303349409654
$ws.Range("T4").Value = 0.3066303349409653

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.288371
$ws.Range("H5").Value = 0.865113
$ws.Range("I5").Value = 0.03053053623570109
$ws.Range("J5").Value = 0.03053053623570109
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.553279333333334
$ws.Range("N5").Value = 7.659838000000001
$ws.Range("O5").Value = 0.1645043904057808
$ws.Range("P5").Value = 0.1645043904057808
$ws.Range("Q5").Value = 0.7362917146326667
$ws.Range("R5").Value = 6.626625431694
$ws.Range("S5").Value = 0.005022407252215609
$ws.Range("T5").Value = 0.005022407252215608

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.288371
$ws.Range("H6").Value = 0.865113
$ws.Range("I6").Value = 0.03053053623570109
$ws.Range("J6").Value = 0.03053053623570109
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.058662
$ws.Range("N6").Value = 24.175986
$ws.Range("O6").Value = 0.5192088709172035
$ws.Range("P6").Value = 0.5192088709172035
$ws.Range("Q6").Value = 2.323884419602
$ws.Range("R6").Value = 20.914959776418
$ws.Range("S6").Value = 0.01585172524743513
$ws.Range("T6").Value = 0.01585172524743513

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.288371
$ws.Range("H7").Value = 0.865113
$ws.Range("I7").Value = 0.03053053623570109
$ws.Range("J7").Value = 0.03053053623570109
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.909099333333334
$ws.Range("N7").Value = 14.727298
$ws.Range("O7").Value = 0.3162867386770157
$ws.Range("P7").Value = 0.3162867386770157
$ws.Range("Q7").Value = 1.415641883852667
$ws.Range("R7").Value = 12.740776954674
$ws.Range("S7").Value = 0.009656403736050348
$ws.Range("T7").Value = 0.009656403736050348
